$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 23 new rows (94:116) BEFORE filling in row 93's Start/End Time, so
#     that the row-format clone (which copies row 93's current layout) does not
#     propagate Start/End Time cells into rows that should not have them.
#     This also picks up the "general" dxf-backed style used by columns E/F. ---
$ws.Range("A94:A116").EntireRow.Insert(-4121)

# --- Fill in row 93 (previously had only a Date; now gets Start/End Time too) ---
$ws.Cells.Item(93,2).Value = 0
$ws.Cells.Item(93,3).Value = 0

# --- Rows 94-97: Date + Start Time (0) + End Time (0), each with its own
#     (non-shared) Duration / Second Duration / Absolute Value formula ---
$dates1 = @(43418, 43419, 43420, 43421)
for ($i = 0; $i -lt $dates1.Length; $i++) {
    $r = 94 + $i
    $ws.Cells.Item($r,1).Value = $dates1[$i]
    $ws.Cells.Item($r,2).Value = 0
    $ws.Cells.Item($r,3).Value = 0
    $ws.Range("D$r").Formula = "=(C$r-B$r)* 1440"
    $ws.Range("E$r").Formula = "=IF(C$r>B$r, (C$r-B$r)*1440, (B$r-C$r)*1440)"
    $ws.Range("F$r").Formula = "=ABS((C$r-B$r)*1440)"
}

# --- Rows 98-100: Date only (no Start/End Time), each with its own
#     (non-shared) formula ---
$dates2 = @(43422, 43423, 43424)
for ($i = 0; $i -lt $dates2.Length; $i++) {
    $r = 98 + $i
    $ws.Cells.Item($r,1).Value = $dates2[$i]
    $ws.Range("D$r").Formula = "=(C$r-B$r)* 1440"
    $ws.Range("E$r").Formula = "=IF(C$r>B$r, (C$r-B$r)*1440, (B$r-C$r)*1440)"
    $ws.Range("F$r").Formula = "=ABS((C$r-B$r)*1440)"
}

# --- Rows 101-104: Date only, sharing one formula per column across the block ---
$dates3 = @(43425, 43426, 43427, 43428)
for ($i = 0; $i -lt $dates3.Length; $i++) {
    $ws.Cells.Item(101 + $i, 1).Value = $dates3[$i]
}
$ws.Range("D101:D104").Formula = "=(C101-B101)* 1440"
$ws.Range("E101:E104").Formula = "=IF(C101>B101, (C101-B101)*1440, (B101-C101)*1440)"
$ws.Range("F101:F104").Formula = "=ABS((C101-B101)*1440)"

# --- Rows 105-110: Date only, sharing one formula per column across the block ---
$dates4 = @(43429, 43430, 43431, 43432, 43433, 43434)
for ($i = 0; $i -lt $dates4.Length; $i++) {
    $ws.Cells.Item(105 + $i, 1).Value = $dates4[$i]
}
$ws.Range("D105:D110").Formula = "=(C105-B105)* 1440"
$ws.Range("E105:E110").Formula = "=IF(C105>B105, (C105-B105)*1440, (B105-C105)*1440)"
$ws.Range("F105:F110").Formula = "=ABS((C105-B105)*1440)"

# --- Rows 111-116: Date only, sharing one formula per column across the block ---
$dates5 = @(43435, 43436, 43437, 43438, 43439, 43440)
for ($i = 0; $i -lt $dates5.Length; $i++) {
    $ws.Cells.Item(111 + $i, 1).Value = $dates5[$i]
}
$ws.Range("D111:D116").Formula = "=(C111-B111)* 1440"
$ws.Range("E111:E116").Formula = "=IF(C111>B111, (C111-B111)*1440, (B111-C111)*1440)"
$ws.Range("F111:F116").Formula = "=ABS((C111-B111)*1440)"

# --- Resize the structured table / autofilter to cover the newly added rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F116"))

# --- Update sheet view (scrolled position / active selection) ---
$ws.Application.ActiveWindow.ScrollRow = 89
$ws.Range("H96").Select()

Write-Output "edit complete"
